$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the header row text with lowercase variants. Because the old
# header strings ("Concejo", "Población", "Zona", "Paisaje_protegido")
# become unreferenced once overwritten, the engine drops them from the
# shared-string table on save and appends the new lowercase strings at
# the end - reproducing the reorder seen in the target workbook.
$ws.Range("A1").Value = "concejo"
$ws.Range("B1").Value = "población"
$ws.Range("C1").Value = "zona"
$ws.Range("D1").Value = "paisaje_protegido"

# Move the active selection to F3, matching the saved view state.
$ws.Range("F3").Select()
